# Update countries & provincias Spain
# Applies the refreshed COVID stats + re-sorted Mauritania row + updated timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Timestamp banner (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 18 de Mayo de 2020 a las 21:35"

# --- Estados Unidos (row 4) ---
$ws.Range("B4").Value = 1542904
$ws.Range("C4").Value = 15240
$ws.Range("D4").Value = 351741
$ws.Range("E4").Value = 1099567
$ws.Range("G4").Value = 618
$ws.Range("H4").Value = 91596

# --- Francia (row 10) ---
$ws.Range("B10").Value = 179927
$ws.Range("C10").Value = 358
$ws.Range("E10").Value = 89960

# --- Alemania (row 11) ---
$ws.Range("B11").Value = 177268
$ws.Range("C11").Value = 617
$ws.Range("E11").Value = 14549
$ws.Range("G11").Value = 70
$ws.Range("H11").Value = 8119

# --- Barein (row 55) ---
$ws.Range("B55").Value = 7184
$ws.Range("C55").Value = 228
$ws.Range("D55").Value = 2931
$ws.Range("E55").Value = 4241

# --- Uzbekistan (row 75) ---
$ws.Range("B75").Value = 2791
$ws.Range("C75").Value = 38
$ws.Range("E75").Value = 464

# --- Sri Lanka (row 105) ---
$ws.Range("B105").Value = 992
$ws.Range("C105").Value = 11
$ws.Range("E105").Value = 424

# --- San Marino (row 122) ---
$ws.Range("D122").Value = 203
$ws.Range("E122").Value = 410

# --- Jordania (row 123) ---
$ws.Range("B123").Value = 629
$ws.Range("C123").Value = 16
$ws.Range("D123").Value = 413
$ws.Range("E123").Value = 207

# --- Rows 172-175: Mauritania's rising case count pushes it above
#     San Martin (Parte Holandesa), Malaui and Libia in the ranking
#     (the sheet is kept sorted by total cases, column B, descending).
#     Mauritania moves up to row 172 with fresh numbers, and the three
#     countries it passes are each pushed down one row, carrying their
#     existing statistics with them. ---

# Row 172: Mauritania (new data)
$ws.Range("A172").Value = "Mauritania"
$ws.Range("B172").Value = 81
$ws.Range("C172").Value = 19
$ws.Range("D172").Value = 7
$ws.Range("E172").Value = 70
$ws.Range("F172").Value = 0
$ws.Range("G172").Value = 0
$ws.Range("H172").Value = 4

# Row 173: San Martin (Parte Holandesa) (previously on row 172)
$ws.Range("A173").Value = "San Martin (Parte Holandesa)"
$ws.Range("B173").Value = 77
$ws.Range("C173").Value = 0
$ws.Range("D173").Value = 54
$ws.Range("E173").Value = 8
$ws.Range("F173").Value = 0
$ws.Range("G173").Value = 0
$ws.Range("H173").Value = 15

# Row 174: Malaui (previously on row 173)
$ws.Range("A174").Value = "Malaui"
$ws.Range("B174").Value = 70
$ws.Range("C174").Value = 0
$ws.Range("D174").Value = 27
$ws.Range("E174").Value = 40
$ws.Range("F174").Value = 0
$ws.Range("G174").Value = 0
$ws.Range("H174").Value = 3

# Row 175: Libia (previously on row 174)
$ws.Range("A175").Value = "Libia"
$ws.Range("B175").Value = 65
$ws.Range("C175").Value = 0
$ws.Range("D175").Value = 35
$ws.Range("E175").Value = 27
$ws.Range("F175").Value = 0
$ws.Range("G175").Value = 0
$ws.Range("H175").Value = 3
